$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (pushes the existing row 6.."谢楠"'s future
# neighbours - formerly rows 6-20, now 7-21 - down by one, and Excel auto
# adjusts the RANK() formula ranges from B$2:B$20 to B$2:B$21).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new contestant.
$ws.Range("A6").Formula = "=RANK(B6,B`$2:B`$21)"
$ws.Range("B6").Value = 61
$ws.Range("C6").Value = "谢楠"

# Two existing contestant names were corrected/renamed. After the row
# insert above, the old row 11 ("zyh") is now row 12, and the old row 13
# ("20200234答案" - a mixed-run rich string) is now row 14. Re-assigning
# Value collapses any rich-text runs to a plain string; force the Arial
# font so the resulting cell style matches the rest of the column
# (the replaced cells previously carried the CJK font style).
$ws.Range("C12").Value = "郑义航"
$ws.Range("C12").Font.Name = "Arial"
$ws.Range("C14").Value = "许婉婷"
$ws.Range("C14").Font.Name = "Arial"

# Restore the view state captured in the saved workbook: zoom level and
# the active cell/selection.
$ws.Application.ActiveWindow.Zoom = 160
$ws.Range("D16").Select()
